{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// precedes them (right after \"LOB1053: F\u00edsica III (Requisito fraco)\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetFooter =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\nconst targetLink = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n\n// Find the \"\u00a9 2020 ...\" paragraph; the two paragraphs immediately before it\n// (the \"Ver no Jupiter...\" line and the blank spacer line) are removed too.\nlet footerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetFooter) {\n    footerIndex = i;\n    break;\n  }\n}\n\nif (footerIndex >= 2 && paragraphs.items[footerIndex - 1].text === targetLink) {\n  paragraphs.items[footerIndex].delete(); // \"\u00a9 2020 ...\" paragraph\n  paragraphs.items[footerIndex - 1].delete(); // \"Ver no Jupiter ...\" paragraph\n  paragraphs.items[footerIndex - 2].delete(); // blank spacer paragraph\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n# precedes them (right after \"LOB1053: F\u00edsica III (Requisito fraco)\").\n$d = $word.ActiveDocument\n\n$targetFooter = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n$targetLink = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$footerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $targetFooter) {\n        $footerIndex = $i\n        break\n    }\n}\n\nif ($footerIndex -ge 3) {\n    $linkText = $d.Paragraphs.Item($footerIndex - 1).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($linkText -eq $targetLink) {\n        $d.Paragraphs.Item($footerIndex).Range.Delete()       # \"\u00a9 2020 ...\" paragraph\n        $d.Paragraphs.Item($footerIndex - 1).Range.Delete()   # \"Ver no Jupiter ...\" paragraph\n        $d.Paragraphs.Item($footerIndex - 2).Range.Delete()   # blank spacer paragraph\n    }\n}\n"}
